$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply cell value updates per the source diff
$ws.Range("D2").Value = '43.280.41'
$ws.Range("E2").Value = '  +5.25%  '
$ws.Range("D3").Value = '2.297.70'
$ws.Range("E3").Value = '  +5.63%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '252.78'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  +5.42%  '
$ws.Range("D7").Value = '73.37'
$ws.Range("E7").Value = '  +11.00%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.666'
$ws.Range("E9").Value = '  +13.60%  '
$ws.Range("D10").Value = '39.81'
$ws.Range("E10").Value = '  +9.82%  '
$ws.Range("E11").Value = '  +5.56%  '
$ws.Range("D12").Value = '59.96'
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '7.55'
$ws.Range("E13").Value = '  +10.34%  '
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = '2.636.99'
$ws.Range("E15").Value = '  +5.45%  '
$ws.Range("D16").Value = '15.21'
$ws.Range("E16").Value = '  +6.94%  '
$ws.Range("D17").Value = '0.902'
$ws.Range("E17").Value = '  +7.11%  '
$ws.Range("D18").Value = '2.288.71'
$ws.Range("E18").Value = '  +5.05%  '
$ws.Range("D19").Value = '43.185.95'
$ws.Range("E19").Value = '  +5.27%  '
$ws.Range("E20").Value = '  +7.63%  '
$ws.Range("D21").Value = '6.43'
$ws.Range("E21").Value = '  +6.69%  '
$ws.Range("D22").Value = '73.83'
$ws.Range("E22").Value = '  +3.35%  '
$ws.Range("D23").Value = '238.68'
$ws.Range("E23").Value = '  +3.96%  '
$ws.Range("D24").Value = '2.19'
$ws.Range("E24").Value = '  +8.06%  '
$ws.Range("B25").Value = 'WEMIXToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D25").Value = '3.92'
$ws.Range("E25").Value = '  +2.39%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '11.90'
$ws.Range("E26").Value = '  +6.11%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").Value = '3.70'
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  +8.57%  '
$ws.Range("D31").Value = '168.35'
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").Value = '21.32'
$ws.Range("E32").Value = '  +5.91%  '
$ws.Range("D33").Value = '6.32'
$ws.Range("E33").Value = '  +12.08%  '
$ws.Range("E34").Value = '  +7.99%  '
$ws.Range("D35").Value = '0.0812'
$ws.Range("E35").Value = '  +9.32%  '
$ws.Range("D36").Value = '31.33'
$ws.Range("E36").Value = '  +28.15%  '
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = '0.127'
$ws.Range("E37").Value = '  +5.34%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '4.79'
$ws.Range("E38").Value = '  +21.06%  '
$ws.Range("D39").Value = '4.83'
$ws.Range("E39").Value = '  +7.32%  '
$ws.Range("E40").Value = '  +3.59%  '
$ws.Range("D41").Value = '13.59'
$ws.Range("E41").Value = '  +21.21%  '
$ws.Range("D42").Value = '2.35'
$ws.Range("E42").Value = '  +7.00%  '
$ws.Range("D43").Value = '6.14'
$ws.Range("E43").Value = '  +11.64%  '
$ws.Range("E44").Value = '  +14.00%  '
$ws.Range("D45").Value = '9.25'
$ws.Range("E45").Value = '  +9.36%  '
$ws.Range("D46").Value = '62.04'
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("D47").Value = '4.96'
$ws.Range("E47").Value = '  -8.77%  '
$ws.Range("E48").Value = '  +5.54%  '
$ws.Range("D49").Value = '1.20'
$ws.Range("E49").Value = '  +6.18%  '
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("E51").Value = '  +5.83%  '
